# ft: removed position cuts
# Updates the "Puntos" standings column on "Equipos" and the match
# Local/Visita teams + Resultado scores on "Resultados".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Equipos")
$ws2 = $wb.Worksheets.Item("Resultados")

# --- Equipos: "Puntos al finalizar la primera rueda" (column D) ---
$ws1.Range("D2").Value = 6
$ws1.Range("D3").Value = 7
$ws1.Range("D4").Value = 4
$ws1.Range("D5").Value = 0

# --- Resultados: swap Local/Visita teams for Jornada 5 matches (rows 6-7) ---
$ws2.Range("C6").Value = "C"
$ws2.Range("D6").Value = "A"
$ws2.Range("C7").Value = "D"
$ws2.Range("D7").Value = "B"

# --- Resultados: updated match scores (column E) ---
$ws2.Range("E3").Value = "2:3"
$ws2.Range("E4").Value = "1:5"
$ws2.Range("E6").Value = "6:4"
$ws2.Range("E7").Value = "0:2"
$ws2.Range("E9").Value = "1:4"
$ws2.Range("E10").Value = "3:2"
$ws2.Range("E12").Value = "0:6"
$ws2.Range("E13").Value = "2:2"
$ws2.Range("E15").Value = "5:4"
$ws2.Range("E16").Value = "1:0"
$ws2.Range("E18").Value = "3:4"
$ws2.Range("E19").Value = "3:1"
